$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.675.90"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "3.456.04"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.71"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.74"
$ws.Range("E6").Value = "  +4.59%  "

$ws.Range("D7").Value = "3.457.24"
$ws.Range("E7").Value = "  +1.58%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("D13").Value = "4.045.22"
$ws.Range("E13").Value = "  +1.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.57"
$ws.Range("E14").Value = "  +7.16%  "

$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "3.458.24"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "61.795.16"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +7.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.45"
$ws.Range("E22").Value = "  +6.47%  "

$ws.Range("E23").Value = "  +2.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.56"
$ws.Range("E24").Value = "  +4.50%  "

$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("E27").Value = "  +1.61%  "

$ws.Range("D28").Value = "3.593.20"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("E29").Value = "  +4.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.66"
$ws.Range("E30").Value = "  +2.92%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.46"
$ws.Range("E34").Value = "  -9.19%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.02"
$ws.Range("E36").Value = "  +2.29%  "

$ws.Range("E37").Value = "  +2.68%  "

$ws.Range("D38").Value = "3.480.02"
$ws.Range("E38").Value = "  +1.76%  "

$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.17"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.83"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0795"
$ws.Range("E42").Value = "  +2.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.20"
$ws.Range("E43").Value = "  +4.26%  "

$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.54"
$ws.Range("E45").Value = "  +3.06%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.45"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").Value = "2.623.42"
$ws.Range("E49").Value = "  +3.62%  "

$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.98"
$ws.Range("E51").Value = "  +2.75%  "

